$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column (2022) in row 4, matching the style of the
# existing "R4" (2021) header cell.
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S4").Value = 2022

# Add the corresponding data value (42) in row 5, matching the style of
# the existing "R5" (42.9) cell, then apply a "0.0" number format to it
# (this introduces the new custom numFmt 166 used only by this cell).
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S5").Value = 42
$ws.Range("S5").NumberFormat = "0.0"

# Match the author's final selection in the sheet view.
[void]$ws.Range("U4").Select()
